$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1718.25
$ws.Range("J17").Value = 1718.25
$ws.Range("L17").Value = 5154.75
$ws.Range("N17").Value = -5490.75

$ws.Range("H53").Value = 504.5625
$ws.Range("I53").Value = 474
$ws.Range("J53").Value = 518.4545000000001
$ws.Range("K53").Value = 474
$ws.Range("L53").Value = 518.4545000000001
$ws.Range("M53").Value = 163
$ws.Range("N53").Value = -1792.4545

$ws.Range("H135").Value = 858.2174
$ws.Range("I135").Value = 865.15
$ws.Range("J135").Value = 812
$ws.Range("K135").Value = 7786.349999999999
$ws.Range("L135").Value = 7308
$ws.Range("M135").Value = -5251.349999999999
$ws.Range("N135").Value = -12378

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4591.8125
$ws.Range("I32").Value = 4060.3428
$ws.Range("J32").Value = 6022.6924
$ws.Range("K32").Value = 4060.3428
$ws.Range("L32").Value = 6022.6924
$ws.Range("M32").Value = -3773.3428
$ws.Range("N32").Value = -6596.6924

$ws.Range("H45").Value = 1429.22
$ws.Range("I45").Value = 977.9143
$ws.Range("J45").Value = 2482.2666
$ws.Range("K45").Value = 977.9143
$ws.Range("L45").Value = 2482.2666
$ws.Range("M45").Value = -600.9143
$ws.Range("N45").Value = -3236.2666

$ws.Range("H61").Value = 4129.278
$ws.Range("I61").Value = 1626.3334
$ws.Range("J61").Value = 4629.8667
$ws.Range("K61").Value = 1626.3334
$ws.Range("L61").Value = 4629.8667
$ws.Range("M61").Value = -1414.3334
$ws.Range("N61").Value = -5053.8667

$ws.Range("H97").Value = 548.6316
$ws.Range("I97").Value = 457.75
$ws.Range("J97").Value = 1033.3334
$ws.Range("K97").Value = 457.75
$ws.Range("L97").Value = 1033.3334
$ws.Range("M97").Value = 38.25
$ws.Range("N97").Value = -2025.3334

$ws.Range("H101").Value = 39280.668
$ws.Range("J101").Value = 39280.668
$ws.Range("L101").Value = 39280.668
$ws.Range("N101").Value = -45770.668

$ws.Range("H109").Value = 32500
$ws.Range("J109").Value = 32500
$ws.Range("L109").Value = 32500
$ws.Range("N109").Value = -35274

$ws.Range("H136").Value = 4129.278
$ws.Range("I136").Value = 1626.3334
$ws.Range("J136").Value = 4629.8667
$ws.Range("K136").Value = 4879.0002
$ws.Range("L136").Value = 13889.6001
$ws.Range("M136").Value = -2329.0002
$ws.Range("N136").Value = -18989.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 558.0714
$ws.Range("J80").Value = 610.4666999999999
$ws.Range("L80").Value = 610.4666999999999
$ws.Range("N80").Value = -2606.4667

$ws.Range("H83").Value = 558.0714
$ws.Range("J83").Value = 610.4666999999999
$ws.Range("L83").Value = 3052.3335
$ws.Range("N83").Value = -13036.3335

$ws.Range("H94").Value = 1637.3055
$ws.Range("I94").Value = 1482.1724
$ws.Range("J94").Value = 2280
$ws.Range("K94").Value = 1482.1724
$ws.Range("L94").Value = 2280
$ws.Range("M94").Value = -1031.1724
$ws.Range("N94").Value = -3182

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2838
$ws.Range("I31").Value = 2198.2144
$ws.Range("J31").Value = 3364.8823
$ws.Range("K31").Value = 2198.2144
$ws.Range("L31").Value = 3364.8823
$ws.Range("M31").Value = -1903.2144
$ws.Range("N31").Value = -3954.8823

$ws.Range("H34").Value = 2838
$ws.Range("I34").Value = 2198.2144
$ws.Range("J34").Value = 3364.8823
$ws.Range("K34").Value = 2198.2144
$ws.Range("L34").Value = 3364.8823
$ws.Range("M34").Value = -1996.2144
$ws.Range("N34").Value = -3768.8823

$ws.Range("H86").Value = 7110.8
$ws.Range("I86").Value = 6000
$ws.Range("J86").Value = 7851.3335
$ws.Range("K86").Value = 6000
$ws.Range("L86").Value = 7851.3335
$ws.Range("M86").Value = -4877
$ws.Range("N86").Value = -10097.3335

$ws.Range("H89").Value = 7110.8
$ws.Range("I89").Value = 6000
$ws.Range("J89").Value = 7851.3335
$ws.Range("K89").Value = 30000
$ws.Range("L89").Value = 39256.6675
$ws.Range("M89").Value = -24384
$ws.Range("N89").Value = -50488.6675

$ws.Range("H134").Value = 3247.05
$ws.Range("I134").Value = 1932.3334
$ws.Range("J134").Value = 4322.727
$ws.Range("K134").Value = 5797.0002
$ws.Range("L134").Value = 12968.181
$ws.Range("M134").Value = -3262.0002
$ws.Range("N134").Value = -18038.181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3749.0625
$ws.Range("J62").Value = 3932.3333
$ws.Range("L62").Value = 11796.9999
$ws.Range("N62").Value = -13168.9999

$ws.Range("H63").Value = 16003.429
$ws.Range("I63").Value = 1006
$ws.Range("J63").Value = 36000
$ws.Range("K63").Value = 3018
$ws.Range("L63").Value = 108000
$ws.Range("M63").Value = -2269
$ws.Range("N63").Value = -109498

$ws.Range("H65").Value = 3749.0625
$ws.Range("J65").Value = 3932.3333
$ws.Range("L65").Value = 35390.9997
$ws.Range("N65").Value = -42254.9997

$ws.Range("H66").Value = 16003.429
$ws.Range("I66").Value = 1006
$ws.Range("J66").Value = 36000
$ws.Range("K66").Value = 9054
$ws.Range("L66").Value = 324000
$ws.Range("M66").Value = -5310
$ws.Range("N66").Value = -331488

$ws.Range("H68").Value = 2222.3157
$ws.Range("I68").Value = 857.5
$ws.Range("J68").Value = 2478.2188
$ws.Range("K68").Value = 2572.5
$ws.Range("L68").Value = 7434.6564
$ws.Range("M68").Value = -1761.5
$ws.Range("N68").Value = -9056.6564

$ws.Range("H71").Value = 2222.3157
$ws.Range("I71").Value = 857.5
$ws.Range("J71").Value = 2478.2188
$ws.Range("K71").Value = 7717.5
$ws.Range("L71").Value = 22303.9692
$ws.Range("M71").Value = -3661.5
$ws.Range("N71").Value = -30415.9692

$ws.Range("H107").Value = 1022.6786
$ws.Range("I107").Value = 629.5333000000001
$ws.Range("J107").Value = 1166.5122
$ws.Range("K107").Value = 1888.5999
$ws.Range("L107").Value = 3499.536599999999
$ws.Range("M107").Value = 31.40009999999984
$ws.Range("N107").Value = -7339.536599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 100003010
$ws.Range("I61").Value = 250003230
$ws.Range("J61").Value = 2868.3333
$ws.Range("K61").Value = 250003230
$ws.Range("L61").Value = 2868.3333
$ws.Range("M61").Value = -250003028
$ws.Range("N61").Value = -3272.3333

$ws.Range("H113").Value = 100003010
$ws.Range("I113").Value = 250003230
$ws.Range("J113").Value = 2868.3333
$ws.Range("K113").Value = 250003230
$ws.Range("L113").Value = 2868.3333
$ws.Range("M113").Value = -250001060
$ws.Range("N113").Value = -7208.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6177.8286
$ws.Range("I132").Value = 2652.84
$ws.Range("K132").Value = 7958.52
$ws.Range("M132").Value = -5428.52

$ws.Range("H136").Value = 3170.5715
$ws.Range("I136").Value = 1476
$ws.Range("J136").Value = 3848.4
$ws.Range("K136").Value = 4428
$ws.Range("L136").Value = 11545.2
$ws.Range("M136").Value = -1878
$ws.Range("N136").Value = -16645.2
